$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 14.45 = 59205.2 pesos`n✅ 59205.2 pesos = 14.39 = 959.5 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update N10, O10, N12, O12 values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 69.2
$ws2.Range("O10").Value = 4097
$ws2.Range("N12").Value = 4115.68
$ws2.Range("O12").Value = 66.7
